# "updates to physio data"
#
# data/physiology/coli_growth_conds.xlsx (Sheet1) edits:
#   1. The sheet had a gap in its row numbering (row 57 never existed -
#      data ran 1..56 then jumped to 58..78). Deleting that blank gap row
#      shifts every fermentation-substrate row up by one (58->57, ..., 77->76).
#   2. Corrects the c_class of trehalose from "sugar" to "sugar alcohol".
#   3. Adds a ref_PMIDS-style citation ("Bergey's manual, Eschericia") in
#      column E for every one of those fermentation rows.
#   4. Replaces the old single "indole / aromatic / oxygen" row with three
#      new aromatic-compound rows (phenylpropanoic acid, phenylacetic acid,
#      3-hydroxy cinnamic acid), all citing PMID 6345502.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the stray empty row 57; everything below slides up one ---
$ws.Rows.Item(57).Delete()

# --- 2. Fix trehalose's c_class (now sitting at row 75) ---
$ws.Range("B75").Value() = "sugar alcohol"

# --- 3. Stamp the Bergey's manual citation onto the fermentation block (rows 57-76) ---
for ($r = 57; $r -le 76; $r++) {
    $ws.Cells.Item($r, 5).Value() = "Bergey's manual, Eschericia"
}

# --- 4. Turn the old single indole/aromatic/oxygen row (now row 77) into
#        three aromatic-compound rows; insert two extra rows to fit them ---
$ws.Rows.Item(78).Insert()
$ws.Rows.Item(78).Insert()

$ws.Range("A77").Value() = "phenylpropanoic acid"
$ws.Range("B77").Value() = "aromatic"
$ws.Range("C77").Value() = "oxygen"
$ws.Range("D77").Value() = 1
$ws.Range("E77").Value() = 6345502

$ws.Range("A78").Value() = "phenylacetic acid"
$ws.Range("B78").Value() = "aromatic"
$ws.Range("C78").Value() = "oxygen"
$ws.Range("D78").Value() = 1
$ws.Range("E78").Value() = 6345502

$ws.Range("A79").Value() = "3-hydroxy cinnamic acid"
$ws.Range("B79").Value() = "aromatic"
$ws.Range("C79").Value() = "oxygen"
$ws.Range("D79").Value() = 1
$ws.Range("E79").Value() = 6345502

# --- 5. Best-effort cosmetic view-state match (scroll position / selection) ---
[void]$ws.Application.Goto($ws.Range("A11"))
[void]$ws.Range("N22").Select()
